{"js": "// Remove the last row of the (only) table in the document \u2014 the\n// \"Bestelt\" row, which documented the old order/product relation that\n// is now obsolete now that \"add new product\" works.\nconst table = context.document.body.tables.getFirst();\nconst lastRow = table.rows.getLast();\nlastRow.delete();\nawait context.sync();\n", "ps1": "# Remove the last row of the (first) table in the document \u2014 the\n# \"Bestelt\" row, which described the old order/product relation that no\n# longer applies now that \"add new product\" works.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$lastRow = $t.Rows.Item($t.Rows.Count)\n$lastRow.Delete()\n"}
